# Add a rotated caption label ("Triple 2-channel analog multiplexer/demultiplexer")
# to the left edge of the logic-diagram picture on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Draw a plain rectangle autoshape (same as using the Shapes > Rectangle tool),
# then size/position/rotate it and turn on "Resize shape to fit text" with no
# wrapping -- matches a one-line auto-fit caption textbox.
$sh = $s.Shapes.AddShape(1, 0, 0, 1, 1)
$sh.Name = "Rectangle 1"

$sh.TextFrame.WordWrap = 0
$sh.TextFrame.AutoSize = 1

$tr = $sh.TextFrame.TextRange
$tr.Text = "Triple 2-channel analog multiplexer/"
$tr.InsertAfter("demultiplexer")

# Unrotated bounding box (points == EMU / 12700), rotated 270 degrees about its
# own center so it reads bottom-to-top along the left edge of the diagram.
$sh.Left = -164.52259842519686
$sh.Top = 226.49503937007873
$sh.Width = 384.0088188976378
$sh.Height = 29.081259842519685
$sh.Rotation = 270
